# Updated symbol list on Mon Jan 16 03:57:57 UTC 2023 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values for the
# cryptocurrency rows that changed, keeping the cells as plain text (the
# workbook stores these numbers/percentages as text, not numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellAddress,
        [string]$Text
    )
    $range = $ws.Range($CellAddress)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

# Row 2 (BNB)
Set-TextValue "D2" "304.19"
Set-TextValue "E2" "1.94%"

# Row 3 (OKB)
Set-TextValue "D3" "31.93"
Set-TextValue "E3" "1.37%"

# Row 4 (HuobiToken)
Set-TextValue "D4" "5.181"
Set-TextValue "E4" "0.29%"

# Row 5 (Cronos)
Set-TextValue "D5" "0.07476"
Set-TextValue "E5" "-0.33%"

# Row 6 (FTXToken)
Set-TextValue "D6" "2.451"
Set-TextValue "E6" "45.71%"

# Row 7 (KuCoinToken)
Set-TextValue "D7" "8.012"
Set-TextValue "E7" "2.96%"

# Row 8 (MXToken)
Set-TextValue "D8" "0.9162"
Set-TextValue "E8" "-0.87%"

# Row 9 (WazirX)
Set-TextValue "D9" "0.1735"
Set-TextValue "E9" "1.24%"

# Row 10 (LiechtensteinCryptoassetsExchange)
Set-TextValue "D10" "0.07690"
Set-TextValue "E10" "2.09%"

# Row 11 (MandalaExchangeToken)
Set-TextValue "D11" "0.08188"
Set-TextValue "E11" "2.74%"

# Row 12 (BitrueCoin)
Set-TextValue "D12" "0.03003"
Set-TextValue "E12" "-1.55%"

# Row 13 (BitMartToken)
Set-TextValue "D13" "0.09925"
Set-TextValue "E13" "0.35%"

# Row 14 (BitForexToken)
Set-TextValue "D14" "0.001510"
Set-TextValue "E14" "0.93%"

# Row 15 (TigerCash)
Set-TextValue "D15" "0.006073"
Set-TextValue "E15" "-6.67%"

# Row 16 (LEO) -- only Volume changes
Set-TextValue "E16" "1.33%"

# Row 17 (GateToken)
Set-TextValue "D17" "3.864"
Set-TextValue "E17" "1.77%"

# Row 18 (BTSEToken)
Set-TextValue "D18" "2.230"
Set-TextValue "E18" "0.02%"

# Row 19 (BitpandaEcosystemToken)
Set-TextValue "D19" "0.3261"
Set-TextValue "E19" "-0.96%"

# Row 20 (ProBitToken) -- only Volume changes
Set-TextValue "E20" "0.29%"

# Row 21 (MCDex)
Set-TextValue "D21" "4.653"
Set-TextValue "E21" "2.15%"

# Row 22 (CoinExToken)
Set-TextValue "D22" "0.04600"
Set-TextValue "E22" "-1.19%"

# Row 23 (ZBToken)
Set-TextValue "D23" "0.1565"
Set-TextValue "E23" "1.03%"

# Row 24 (BitKan)
Set-TextValue "D24" "0.001262"
Set-TextValue "E24" "3.46%"

# Row 25 (HotbitToken)
Set-TextValue "D25" "0.004523"
Set-TextValue "E25" "2.43%"

# Row 26 (NitroEx)
Set-TextValue "D26" "0.0001299"
Set-TextValue "E26" "-7.14%"

# Row 27 (UpBots)
Set-TextValue "D27" "0.0002741"
Set-TextValue "E27" "51.59%"

# Row 39 (One)
Set-TextValue "D39" "0.01765"
Set-TextValue "E39" "6.41%"

# Row 40 (IDEX)
Set-TextValue "D40" "0.04540"
Set-TextValue "E40" "0.17%"

# Row 41 (KickToken)
Set-TextValue "D41" "0.007398"
Set-TextValue "E41" "7.04%"

# Row 42 (BKEXToken) -- only Volume changes
Set-TextValue "E42" "1.46%"

# Row 43 (CEJI)
Set-TextValue "D43" "0.002129"
Set-TextValue "E43" "3.41%"

# Row 44 (LocalTraders)
Set-TextValue "D44" "0.01078"
Set-TextValue "E44" "-16.01%"

# Row 45 (CoinLion)
Set-TextValue "D45" "0.00006450"
Set-TextValue "E45" "5.93%"

# Row 46 (BOLO) -- only Volume changes
Set-TextValue "E46" "15.26%"
